# Edit script: rename stat sheets to friendlier titles, bump the "Age"
# column (years-days) by one day across all player tables, and fix the
# StandardStats / PlayingTime sheet's "Playing Time" merged header so the
# label sits over columns G:I instead of F:I (with a new placeholder label
# in F1).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename worksheets (Matches and Possession keep their names).
# ---------------------------------------------------------------------
$renames = @{
    "StandardStats"    = "Standard Stats";
    "ShootingStats"    = "Shooting Stats";
    "PassingStats"     = "Passing Stats";
    "PassTypes"        = "Pass Types";
    "GoalShotCreation" = "Goal & Shot Creation";
    "DefensiveActions" = "Defensive Actions";
    "PlayingTime"      = "Playing Time";
    "MiscStats"        = "Miscellaneous Stats";
}

foreach ($oldName in $renames.Keys) {
    $ws = $wb.Worksheets.Item($oldName)
    $ws.Name = $renames[$oldName]
}

# ---------------------------------------------------------------------
# 2) Bump every "Age" cell (format "YY-DDD") in column E by one day on
#    every player-stats sheet (everything except "Matches").
# ---------------------------------------------------------------------
$statSheetNames = @(
    "Standard Stats",
    "Shooting Stats",
    "Passing Stats",
    "Pass Types",
    "Goal & Shot Creation",
    "Defensive Actions",
    "Possession",
    "Playing Time",
    "Miscellaneous Stats"
)

foreach ($sheetName in $statSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($r = 4; $r -le 40; $r++) {
        $cell = $ws.Range("E$r")
        $v = $cell.Value2
        if ($v -match '^(\d+)-(\d+)$') {
            $years = $matches[1]
            $days = [int]$matches[2] + 1
            $newVal = "$years-" + $days.ToString().PadLeft(3, '0')
            $cell.Value = $newVal
        }
    }
}

# ---------------------------------------------------------------------
# 3) On "Standard Stats" and "Playing Time" sheets, the "Playing Time"
#    group header currently spans F1:I1. Shift it to G1:I1 and give F1
#    its own (previously-implicit) "Unnamed: 4_level_0" label.
# ---------------------------------------------------------------------
$headerFixSheets = @("Standard Stats", "Playing Time")

foreach ($sheetName in $headerFixSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F1:I1").UnMerge()
    $ws.Range("F1").Value = "Unnamed: 4_level_0"
    $ws.Range("G1").Value = "Playing Time"
    $ws.Range("G1:I1").Merge()
    # Re-normalize the border so all four cells keep the original header
    # style instead of Excel's automatic merge/unmerge border splitting.
    $ws.Range("F1:I1").Borders.LineStyle = 1
}
